$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($rng, $newName) {
    $count = $rng.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShp = $rng.InlineShapes.Item($i)
        $floatShp = $inlineShp.ConvertToShape()
        $floatShp.Name = $newName
        [void]$floatShp.ConvertToInlineShape()
    }
}

# Footers: primary (Item 1) and first-page (Item 2) both carry the Pearson
# logo picture currently named "image1.png" -> rename to "image2.png".
Rename-InlineShape $sec.Footers.Item(1).Range "image2.png"
Rename-InlineShape $sec.Footers.Item(2).Range "image2.png"

# Header (first-page, Item 2) carries the BTEC logo picture currently named
# "image2.jpg" -> rename to "image1.jpg".
Rename-InlineShape $sec.Headers.Item(2).Range "image1.jpg"
